$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header date in B1 (leading apostrophe forces text so Excel doesn't
# reinterpret the dd/mm/yyyy string as a date serial value)
$ws.Range("B1").Value = "'07/04/2023"

# Remove the "COV-NOTTE" row (row 16). This shifts rows 17-20 up to 16-19.
$ws.Rows.Item(16).Delete()

# Update B/C values for rows 2-15 (fcst values changed, hour moved 14 -> 16)
$ws.Range("B2").Value = 1159
$ws.Range("C2").Value = 16

$ws.Range("B3").Value = 341
$ws.Range("C3").Value = 16

$ws.Range("B4").Value = 15
$ws.Range("C4").Value = 16

$ws.Range("B5").Value = 294
$ws.Range("C5").Value = 16

$ws.Range("B6").Value = 70
$ws.Range("C6").Value = 16

$ws.Range("B7").Value = 111
$ws.Range("C7").Value = 16

$ws.Range("B8").Value = 407
$ws.Range("C8").Value = 16

$ws.Range("B9").Value = 54
$ws.Range("C9").Value = 16

$ws.Range("B10").Value = 815
$ws.Range("C10").Value = 16

$ws.Range("B11").Value = 376
$ws.Range("C11").Value = 16

$ws.Range("B12").Value = 980
$ws.Range("C12").Value = 16

$ws.Range("B13").Value = 564
$ws.Range("C13").Value = 16

$ws.Range("B14").Value = 926
$ws.Range("C14").Value = 16

$ws.Range("B15").Value = 274
$ws.Range("C15").Value = 16

# After the row-16 deletion, rows 16-19 now hold what used to be rows 17-20.
# Update their B/C values as per new data (labels are already correct after shift).
$ws.Range("B16").Value = 141
$ws.Range("C16").Value = 16

$ws.Range("B17").Value = 4
$ws.Range("C17").Value = 16

$ws.Range("B18").Value = 19
$ws.Range("C18").Value = 16

$ws.Range("B19").Value = 77
$ws.Range("C19").Value = 16
